# Add a new "price snapshot" column to the LDLC price-history sheet.
#
# The sheet tracks one timestamped price column per scrape run, followed
# by two trailing columns "nom" (product name) and "url_produit" (product
# URL). A new scrape ("2026-01-29 08:25:15") needs to be inserted just
# before those trailing columns, i.e. a new column is inserted at AK,
# pushing the old AK ("nom") to AL and the old AL ("url_produit") to AM.
#
# For the header row, the new AK1 cell gets the new timestamp string.
# For every data row, the new AK cell gets a copy of the previous last
# price column (AJ) - i.e. the most recent known price is carried
# forward into the new snapshot column (same value, and same "no data"
# blank state, as column AJ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column in front of the existing "AK" column (shifts the
# old AK/AL -> AL/AM, exactly like Excel's own "Insert Column" command).
$ws.Range("AK1").EntireColumn.Insert()

# Header for the freshly inserted column: the new scrape's timestamp.
$ws.Range("AK1").Value = "2026-01-29 08:25:15"

# Carry the latest known price (column AJ, the scrape right before this
# one) into the new column AK for every data row.
$lastRow = 206
for ($r = 2; $r -le $lastRow; $r++) {
    $previousPrice = $ws.Cells.Item($r, 36).Value()
    $ws.Cells.Item($r, 37).Value = $previousPrice
}
